$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws3 = $wb.Worksheets.Item("本地生活")
$ws4 = $wb.Worksheets.Item("全部类型")

$ws1.Range("F3").Value = 163
$ws1.Range("F4").Value = 1791
$ws1.Range("F5").Value = 3342
$ws1.Range("F6").Value = 1126
$ws1.Range("F7").Value = 2245
$ws1.Range("F8").Value = 2165
$ws1.Range("F9").Value = 1133
$ws1.Range("F12").Value = 1694
$ws1.Range("F16").Value = 309
$ws1.Range("F18").Value = 1607
$ws1.Range("F19").Value = 264
$ws1.Range("F20").Value = 859
$ws1.Range("F22").Value = 267
$ws1.Range("F23").Value = 627
$ws1.Range("F24").Value = 12367
$ws1.Range("F25").Value = 12423
$ws1.Range("F27").Value = 713
$ws1.Range("F29").Value = 258
$ws1.Range("F31").Value = 403
$ws1.Range("G33").Value = 50
$ws1.Range("F34").Value = 8
$ws2.Range("F6").Value = 120
$ws3.Range("F3").Value = 97
$ws4.Range("F4").Value = 163
$ws4.Range("F5").Value = 1791
$ws4.Range("F6").Value = 3342
$ws4.Range("F7").Value = 1126
$ws4.Range("F8").Value = 2245
$ws4.Range("F9").Value = 2165
$ws4.Range("F10").Value = 1133
$ws4.Range("F12").Value = 97
$ws4.Range("F14").Value = 1694
$ws4.Range("F20").Value = 309
$ws4.Range("F23").Value = 1607
$ws4.Range("F24").Value = 264
$ws4.Range("F25").Value = 859
$ws4.Range("F27").Value = 267
$ws4.Range("F28").Value = 627
$ws4.Range("F29").Value = 12367
$ws4.Range("F30").Value = 12423
$ws4.Range("F34").Value = 258
$ws4.Range("F36").Value = 403
$ws4.Range("G39").Value = 50
$ws4.Range("F40").Value = 120
$ws4.Range("F41").Value = 8
